# Updates the "cryptos" price/volume table to the latest scraped snapshot.
# Mirrors the GitHub Actions commit "Updated cryptos list ... with GitHub Actions".
#
# Columns: A=index, B=Coin, C=Link, D=Price (text), E=Volume(1h) (text, padded w/ spaces).
# Note: several "Price" values are plain-looking decimals (e.g. 241.04). The sheet
# stores these as TEXT, so we force NumberFormat "@" (Text) before assigning those
# specific cells to stop Excel from auto-coercing them into numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "97.286.71"
$ws.Range("E2").Value = "  +2.27%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.578.63"
$ws.Range("E3").Value = "  +0.11%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 - Solana
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.04"
$ws.Range("E5").Value = "  +2.33%  "

# Row 6 - was XRP, now BNB (rows 6/7 swapped rank order)
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "655.28"
$ws.Range("E6").Value = "  -0.40%  "

# Row 7 - was BNB, now XRP
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.72"
$ws.Range("E7").Value = "  +16.76%  "

# Row 8 - Dogecoin
$ws.Range("E8").Value = "  +7.89%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  -0.08%  "

# Row 10 - Cardano
$ws.Range("E10").Value = "  +4.54%  "

# Row 11 - LidoStakedEther
$ws.Range("D11").Value = "3.575.25"
$ws.Range("E11").Value = "  +0.02%  "

# Row 12 - Avalanche
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "44.29"
$ws.Range("E12").Value = "  +4.36%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +0.45%  "

# Row 14 - Toncoin
$ws.Range("E14").Value = "  -0.45%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("E15").Value = "  +0.15%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "97.007.02"
$ws.Range("E16").Value = "  +2.10%  "

# Row 18 - Polkadot
$ws.Range("E18").Value = "  +11.93%  "

# Row 19 - WrappedEther
$ws.Range("D19").Value = "3.580.08"
$ws.Range("E19").Value = "  +0.27%  "

# Row 20 - Uniswap
$ws.Range("E20").Value = "  +0.22%  "

# Row 21 - Chainlink
$ws.Range("E21").Value = "  +1.15%  "

# Row 22 - Stellar
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.533"
$ws.Range("E22").Value = "  +11.30%  "

# Row 23 - SuiNetwork
$ws.Range("E23").Value = "  +0.84%  "

# Row 24 - BitcoinCash
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "514.66"
$ws.Range("E24").Value = "  +1.13%  "

# Row 25 - PEPE
$ws.Range("E25").Value = "  +5.38%  "

# Row 26 - NEARProtocol
$ws.Range("E26").Value = "  +0.45%  "

# Row 27 - Litecoin
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "101.82"
$ws.Range("E27").Value = "  +6.99%  "

# Row 28 - Aptos
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "13.10"
$ws.Range("E28").Value = "  +3.36%  "

# Row 29 - WrappedeETH
$ws.Range("D29").Value = "3.771.03"
$ws.Range("E29").Value = "  +0.16%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  +17.44%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -1.59%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("E32").Value = "  +3.92%  "

# Row 33 - Dai
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.05%  "

# Row 34 - Cronos
$ws.Range("E34").Value = "  +4.49%  "

# Row 35 - Binance-PegBSC-USD
$ws.Range("E35").Value = "  +0.12%  "

# Row 36 - EthereumClassic
$ws.Range("E36").Value = "  -0.04%  "

# Row 37 - RenderToken
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.78"
$ws.Range("E37").Value = "  +3.86%  "

# Row 38 - Bittensor
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "613.28"
$ws.Range("E38").Value = "  +5.61%  "

# Row 39 - PolygonEcosystemToken
$ws.Range("E39").Value = "  +1.51%  "

# Row 40 - Fetch.AI
$ws.Range("E40").Value = "  -1.79%  "

# Row 41 - Kaspa
$ws.Range("E41").Value = "  +2.60%  "

# Row 42 - ImmutableX
$ws.Range("E42").Value = "  +5.77%  "

# Row 43 - USDe
$ws.Range("E43").Value = "  -0.09%  "

# Row 44 - ARBITRUM
$ws.Range("E44").Value = "  +1.98%  "

# Row 45 - Filecoin
$ws.Range("E45").Value = "  +4.76%  "

# Row 46 - VeChain
$ws.Range("E46").Value = "  +6.14%  "

# Row 47 - Stacks
$ws.Range("E47").Value = "  +2.07%  "

# Row 48 - Algorand
$ws.Range("E48").Value = "  +36.55%  "

# Row 49 - WhiteBITCoin
$ws.Range("E49").Value = "  +0.96%  "

# Row 50 - Cosmos
$ws.Range("E50").Value = "  +4.37%  "

# Row 51 - dogwifhat
$ws.Range("E51").Value = "  +7.53%  "
